$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("A9").Value = "allianz"
$ws.Range("B9").Value = 61
$ws.Range("C9").Value = 42465

$ws.Range("C8").Copy()
$ws.Range("C9").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

$ws.Range("C10").Select()
